$wb = $excel.ActiveWorkbook

# Rename the header cell (and thus the table column name) on each sheet
# from the sheet-specific "<x>_names" label to the common "names" label.
$wsPollutant = $wb.Worksheets.Item("pollutant")
$wsPollutant.Range("A1").Value = "names"

$wsNfr = $wb.Worksheets.Item("nfr")
$wsNfr.Range("A1").Value = "names"

$wsGnfr = $wb.Worksheets.Item("gnfr")
$wsGnfr.Range("A1").Value = "names"

# The file was last saved with the "pollutant" tab active.
$wsPollutant.Activate()
